$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.103.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.11%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3032"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07428"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("E10").Value = "  -6.60%  "

$ws.Range("E11").Value = "  -2.86%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7244"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.36%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.218"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.099.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.790"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007665"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.088.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.570"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1472"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.940"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("E28").Value = "  -3.72%  "

$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("E30").Value = "  -8.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.463"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.489"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.009"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05197"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.32%  "

$ws.Range("E35").Value = "  -5.57%  "

$ws.Range("E36").Value = "  -6.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.646"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01872"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.670"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9058"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4292"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.907"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.048.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.107"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.84%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.183"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.07%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.980.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.82%  "

